$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Authors" column (E) for rows 2-14 holds references imported from an
# elastic-search JSON export. The new import pass re-joins the author list
# with two extra spaces inserted after every comma separator compared to
# the value currently stored in the workbook. Re-derive each new value from
# the existing one instead of hard-coding the (very long) strings.
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value()
    $updated = $current -replace ',(\s+)', ',$1  '
    $cell.Value = $updated
}
